$wb = $excel.ActiveWorkbook

# New data point for 6 April 2020 (serial date 43927) added to all three
# sheets. Each gets row 31: the date, a running-total formula in B, and
# the day's raw count in C - the same pattern the existing rows use.
# Formatting is carried forward from row 30 via copy/paste-format.

# ---- Confirmed sheet ----
$ws1 = $wb.Worksheets.Item("Confirmed")
$ws1.Range("A30:C30").Copy() | Out-Null
$ws1.Range("A31:C31").PasteSpecial(-4122) | Out-Null
$ws1.Range("A31").Value = 43927
$ws1.Range("B31").Formula = "=SUM(B30+C31)"
$ws1.Range("C31").Value = 35
$ws1.Range("B31").Select() | Out-Null

# ---- Recoverd sheet ----
$ws2 = $wb.Worksheets.Item("Recoverd")
$ws2.Range("A30:C30").Copy() | Out-Null
$ws2.Range("A31:C31").PasteSpecial(-4122) | Out-Null
$ws2.Range("A31").Value = 43927
$ws2.Range("B31").Formula = "=SUM(B30+C31)"
$ws2.Range("C31").Value = 3
$ws2.Range("B31").Select() | Out-Null

# ---- Death sheet ----
# (column B here normally carries the "center" style s=1, but the new
# row's B cell instead picks up the "center+middle" style s=2 that
# columns A/C already use, so pull that format from C30 instead.)
$ws3 = $wb.Worksheets.Item("Death")
$ws3.Range("A30:C30").Copy() | Out-Null
$ws3.Range("A31:C31").PasteSpecial(-4122) | Out-Null
$ws3.Range("C30").Copy() | Out-Null
$ws3.Range("B31").PasteSpecial(-4122) | Out-Null
$ws3.Range("A31").Value = 43927
$ws3.Range("B31").Formula = "=SUM(B30+C31)"
$ws3.Range("C31").Value = 3
$ws3.Range("B34").Select() | Out-Null
